# Weekly update: insert two new "Repollo" price rows (Crespo record,
# Primera/Segunda) at the top of the data block, pushing the existing
# rows down by two. This mirrors the upstream weekly refresh where the
# newest observations are prepended to the historical series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row of the block
# (row 432). Excel shifts rows 432:526 down to 434:528, carrying their
# content and formatting with them (dimension grows to R528 automatically).
$ws.Rows("432:433").Insert()

# --- New row 432 (Crespo record / Primera) ---
$ws.Range("A432").Value = 7
$ws.Range("B432").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C432").Value = "Ñuble"
$ws.Range("D432").Value = 45211
$ws.Range("E432").Value = 16
$ws.Range("F432").Value = 100112006
$ws.Range("G432").Value = "Repollo"
$ws.Range("H432").Value = "Crespo record"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 800
$ws.Range("K432").Value = 900
$ws.Range("L432").Value = 1000
$ws.Range("M432").Value = 962
$ws.Range("N432").Value = "$/unidad"
$ws.Range("O432").Value = "Provincia de Diguillín"
$ws.Range("P432").Value = 962
$ws.Range("Q432").Value = 1
$ws.Range("R432").Value = "Hortaliza"

# --- New row 433 (Crespo record / Segunda) ---
$ws.Range("A433").Value = 7
$ws.Range("B433").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C433").Value = "Ñuble"
$ws.Range("D433").Value = 45211
$ws.Range("E433").Value = 16
$ws.Range("F433").Value = 100112006
$ws.Range("G433").Value = "Repollo"
$ws.Range("H433").Value = "Crespo record"
$ws.Range("I433").Value = "Segunda"
$ws.Range("J433").Value = 600
$ws.Range("K433").Value = 700
$ws.Range("L433").Value = 800
$ws.Range("M433").Value = 750
$ws.Range("N433").Value = "$/unidad"
$ws.Range("O433").Value = "Provincia de Diguillín"
$ws.Range("P433").Value = 750
$ws.Range("Q433").Value = 1
$ws.Range("R433").Value = "Hortaliza"
